{"js": "// Apply the \"Added many more features\" edits to the Aloha Fruit Bonanza\n// review document:\n//   1. Title (Heading 1) and the closing bold \"title\" run:\n//        \"Play Aloha Fruit Bonanza Free - Review & Guide | RTP 97.01%\"\n//        -> \"Play Aloha Fruit Bonanza for Free\"\n//   2. \"What we like\" bullet:\n//        \"Cascading wins create new winning combos\"\n//        -> \"Cascading wins create new winning combinations\"\n//   3. \"What we like\" bullet:\n//        \"Stunning design and visuals to create a tropical paradise\"\n//        -> \"Stunning design and visuals\"\n//   4. \"What we don't like\" bullet:\n//        \"Juicy Multipliers are not available during free spins\"\n//        -> \"Juicy Multipliers not available during free spins\"\n//   5. \"What we don't like\" bullet:\n//        \"Limited number of free spins available\"\n//        -> \"Limited number of free spins\"\n//   6. Closing italic \"meta description\" run:\n//        \"Discover the stunning visuals and innovative Scatter Pays feature\n//         in Aloha Fruit Bonanza. With an RTP of 97.01%, play for free and\n//         win up to 7,500x your stake.\"\n//        -> \"Read our review of Aloha Fruit Bonanza and play for free.\n//            Discover exciting gameplay features and stunning design.\"\n\nconst body = context.document.body;\n\n// Pairs of [oldText, newText]. The title/bold pair appears twice in the\n// document (the Heading 1 and the bold run near the end) and search()\n// will return both occurrences, which we replace in a single pass.\nconst replacements = [\n  [\n    \"Play Aloha Fruit Bonanza Free - Review & Guide | RTP 97.01%\",\n    \"Play Aloha Fruit Bonanza for Free\",\n  ],\n  [\n    \"Cascading wins create new winning combos\",\n    \"Cascading wins create new winning combinations\",\n  ],\n  [\n    \"Stunning design and visuals to create a tropical paradise\",\n    \"Stunning design and visuals\",\n  ],\n  [\n    \"Juicy Multipliers are not available during free spins\",\n    \"Juicy Multipliers not available during free spins\",\n  ],\n  [\n    \"Limited number of free spins available\",\n    \"Limited number of free spins\",\n  ],\n  [\n    \"Discover the stunning visuals and innovative Scatter Pays feature in Aloha Fruit Bonanza. With an RTP of 97.01%, play for free and win up to 7,500x your stake.\",\n    \"Read our review of Aloha Fruit Bonanza and play for free. Discover exciting gameplay features and stunning design.\",\n  ],\n];\n\n// Run every search first so all `load()` calls can be satisfied by a single\n// sync, then perform the replacements in a second pass.\nconst searchResults = replacements.map(([oldText]) =>\n  body.search(oldText, { matchCase: true, matchWholeWord: false })\n);\nsearchResults.forEach((results) => results.load(\"items\"));\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newText] = replacements[i];\n  const items = searchResults[i].items;\n  for (let j = 0; j < items.length; j++) {\n    items[j].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Apply the \"Added many more features\" edits to the Aloha Fruit Bonanza\n# review document:\n#   1. Title (Heading 1) and the closing bold \"title\" run:\n#        \"Play Aloha Fruit Bonanza Free - Review & Guide | RTP 97.01%\"\n#        -> \"Play Aloha Fruit Bonanza for Free\"\n#   2. \"What we like\" bullet:\n#        \"Cascading wins create new winning combos\"\n#        -> \"Cascading wins create new winning combinations\"\n#   3. \"What we like\" bullet:\n#        \"Stunning design and visuals to create a tropical paradise\"\n#        -> \"Stunning design and visuals\"\n#   4. \"What we don't like\" bullet:\n#        \"Juicy Multipliers are not available during free spins\"\n#        -> \"Juicy Multipliers not available during free spins\"\n#   5. \"What we don't like\" bullet:\n#        \"Limited number of free spins available\"\n#        -> \"Limited number of free spins\"\n#   6. Closing italic \"meta description\" run:\n#        \"Discover the stunning visuals and innovative Scatter Pays feature\n#         in Aloha Fruit Bonanza. With an RTP of 97.01%, play for free and\n#         win up to 7,500x your stake.\"\n#        -> \"Read our review of Aloha Fruit Bonanza and play for free.\n#            Discover exciting gameplay features and stunning design.\"\n\n$d = $word.ActiveDocument\n\n# wdReplaceAll swaps every occurrence of Find.Text for Find.Replacement.Text\n# (the title/bold pair occurs twice - in the Heading 1 and in the bold run\n# near the end of the document - and a single Execute call with wdReplaceAll\n# handles both).\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nfunction Replace-AllText($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, [ref]$find.Replacement.Text, $wdReplaceAll)\n}\n\nReplace-AllText \"Play Aloha Fruit Bonanza Free - Review & Guide | RTP 97.01%\" \"Play Aloha Fruit Bonanza for Free\"\nReplace-AllText \"Cascading wins create new winning combos\" \"Cascading wins create new winning combinations\"\nReplace-AllText \"Stunning design and visuals to create a tropical paradise\" \"Stunning design and visuals\"\nReplace-AllText \"Juicy Multipliers are not available during free spins\" \"Juicy Multipliers not available during free spins\"\nReplace-AllText \"Limited number of free spins available\" \"Limited number of free spins\"\nReplace-AllText \"Discover the stunning visuals and innovative Scatter Pays feature in Aloha Fruit Bonanza. With an RTP of 97.01%, play for free and win up to 7,500x your stake.\" \"Read our review of Aloha Fruit Bonanza and play for free. Discover exciting gameplay features and stunning design.\"\n"}
